# Updates the "cryptos" price/volume snapshot (Coin, Link, Price, Volume(1h))
# to the latest scraped values. Price-column cells that look like plain
# decimals (e.g. "555.18") are written with a leading apostrophe so Excel
# keeps them as text instead of auto-converting to Number, then the style
# is reset to "Normal" so no stray text-format style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.844.22"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "3.337.99"
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'555.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'174.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  -2.63%  "
$ws.Range("D8").Value = "3.330.17"
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("D12").Value = "'53.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'8.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "3.871.82"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "'18.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "3.350.90"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "63.767.28"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'0.973"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "'431.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("E23").Value = "  +10.31%  "
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'83.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "'13.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("D27").Value = "'10.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "'2.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "'8.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'29.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").Value = "'590.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "'11.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "'58.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("D38").Value = "'3.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("D39").Value = "'35.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("D42").Value = "3.114.93"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'3.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "'0.0405"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'133.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
